# Edit script for 3-5-timer.docx
# Implements the text changes described in the commit:
#  - Removed "specific" before "function" (it runs your specific function -> it runs your function)
#  - Added comma after "Specifically"
#  - Added comma after "a timer" in "To create a timer you call"
#  - Added Oxford comma after "the function you want called"
#  - Inserted "in " into "For instance, exercise 5" -> "For instance, in exercise 5"
#  - Removed comma after "reading" in "for instance reading, an I2C"
#  - Reworded "these timers are CLOSE ... because they won't run until it is their turn"
#    to "the time is CLOSE ... because the timer function won't run until its turn"

$d = $word.ActiveDocument
$wdReplaceAll = 2

function Replace-Text($find, $replace) {
    $r = $d.Content
    $ok = $r.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, $wdReplaceAll)
    if (-not $ok) {
        throw "Find/Replace failed for: $find"
    }
}

# 1. "it runs your specific function" -> "it runs your function"
Replace-Text "it runs your specific function" "it runs your function"

# 2. "does.  Specifically it runs" -> "does.  Specifically, it runs"
Replace-Text "does.  Specifically it runs" "does.  Specifically, it runs"

# 3. "To create a timer you call" -> "To create a timer, you call"
Replace-Text "To create a timer you call" "To create a timer, you call"

# 4. "the function you want called and the argument" -> "the function you want called, and the argument"
Replace-Text "the function you want called and the argument" "the function you want called, and the argument"

# 5. "For instance, exercise 5" -> "For instance, in exercise 5"
Replace-Text "For instance, exercise 5" "For instance, in exercise 5"

# 6. "for instance reading, an I2C" -> "for instance reading an I2C"
Replace-Text "for instance reading, an I2C" "for instance reading an I2C"

# 7. Reword the RTOS timer-accuracy sentence.
Replace-Text "these timers are CLOSE to correct but are not perfectly correct because they won" "the time is CLOSE to correct but is not perfectly correct because the timer function won"

# 8. "won't run until it is their turn" -> "won't run until its turn"
Replace-Text "run until it is their turn" "run until its turn"

# 9. Move the _GoBack bookmark to follow the now-last-edited paragraph (after "...regular basis.  ")
$d.Bookmarks.Item("_GoBack").Delete()
$r = $d.Content
$found = $r.Find.Execute("some regular basis.  ", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $gobackRange = $d.Range($r.End, $r.End)
    $d.Bookmarks.Add("_GoBack", $gobackRange)
}
